$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

$ws.Cells.Item($row, 1).Value = "2024-02-14_23-21-55"
$ws.Cells.Item($row, 2).Value = "yolov5"
$ws.Cells.Item($row, 3).Value = 8
$ws.Cells.Item($row, 4).Value = 1024
$ws.Cells.Item($row, 5).Value = 2400
$ws.Cells.Item($row, 6).Value = 600
$ws.Cells.Item($row, 7).Value = 10
$ws.Cells.Item($row, 8).Value = "'"
$ws.Cells.Item($row, 9).Value = 0.018909
$ws.Cells.Item($row, 10).Value = 0.026391
$ws.Cells.Item($row, 11).Value = 0.0054985
$ws.Cells.Item($row, 12).Value = 0.53124
$ws.Cells.Item($row, 13).Value = 0.77271
$ws.Cells.Item($row, 14).Value = 0.6584
$ws.Cells.Item($row, 15).Value = 0.33023
$ws.Cells.Item($row, 16).Value = 10
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = 522.0679485999999
$ws.Cells.Item($row, 19).Value = 30
$ws.Cells.Item($row, 20).Value = 50
$ws.Cells.Item($row, 21).Value = 3072
$ws.Cells.Item($row, 22).Value = 0.0000001
$ws.Cells.Item($row, 23).Value = 6.5
$ws.Cells.Item($row, 24).Value = 131.4
$ws.Cells.Item($row, 25).Value = 0.00016
$ws.Cells.Item($row, 26).Value = 0.0000001019271325721766
$ws.Cells.Item($row, 27).Value = 3072
